# Update "想去人数" (want-to-go count) figures across the four sheets of the
# workbook to reflect the latest scrape snapshot (gh-pages output @ 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 7543
$ws1.Range("F3").Value  = 91
$ws1.Range("F5").Value  = 4438
$ws1.Range("F8").Value  = 605
$ws1.Range("F9").Value  = 440
$ws1.Range("F10").Value = 132
$ws1.Range("F11").Value = 342
$ws1.Range("F12").Value = 753
$ws1.Range("F14").Value = 61
$ws1.Range("F15").Value = 238
$ws1.Range("F17").Value = 242
$ws1.Range("F19").Value = 377
$ws1.Range("F20").Value = 137
$ws1.Range("F21").Value = 1079
$ws1.Range("F23").Value = 460
$ws1.Range("F24").Value = 2145
$ws1.Range("F25").Value = 671
$ws1.Range("F26").Value = 32
$ws1.Range("F27").Value = 33
$ws1.Range("F28").Value = 37
$ws1.Range("F29").Value = 590
$ws1.Range("F30").Value = 37

# --- Sheet: 演出 ---------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 284

# --- Sheet: 本地生活 -------------------------------------------------------
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 424

# --- Sheet: 全部类型 -------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 424
$ws4.Range("F3").Value  = 7543
$ws4.Range("F4").Value  = 91
$ws4.Range("F6").Value  = 284
$ws4.Range("F7").Value  = 4438
$ws4.Range("F10").Value = 605
$ws4.Range("F11").Value = 440
$ws4.Range("F13").Value = 132
$ws4.Range("F14").Value = 342
$ws4.Range("F18").Value = 753
$ws4.Range("F20").Value = 62
$ws4.Range("F21").Value = 238
$ws4.Range("F26").Value = 242
$ws4.Range("F28").Value = 377
$ws4.Range("F29").Value = 137
$ws4.Range("F30").Value = 1079
$ws4.Range("F32").Value = 460
$ws4.Range("F33").Value = 2145
$ws4.Range("F34").Value = 671
$ws4.Range("F35").Value = 32
$ws4.Range("F36").Value = 33
$ws4.Range("F37").Value = 37
$ws4.Range("F38").Value = 590
$ws4.Range("F39").Value = 37
